$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "TestAuto_POC1"
$ws.Range("B5").Value = "TestAuto_POC1"
$ws.Range("C5").Value = "Facility_POC1"
$ws.Range("D5").Value = "Facility_POC1"
$ws.Range("E5").Value = "Pharmacy_POC1"
$ws.Range("F5").Value = "Pharmacy_POC1"
$ws.Range("H5").Value = "AlignmentProject_POC1"

$ws.Range("H5").Select()
